$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Edit 1: merge the two runs "MON Dec 18" + " 11:19:08 PST 2017"
# into a single run "MON Dec 18 11:19:08 PST 2017".
# -----------------------------------------------------------------
$found = $d.Content.Find.Execute(
    "MON Dec 18` 11:19:08 PST 2017",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "MON Dec 18 11:19:08 PST 2017", 2)
Write-Host "Edit1 (merge date run) found/replaced: $found"

# -----------------------------------------------------------------
# Edit 2: after the final "Amount Received mode ... - CASH" entry,
# add a brand new purchase-details block (08/02/2018 check-in).
# -----------------------------------------------------------------

# Locate the last paragraph whose text is exactly
# "Amount Received mode" <tab><tab> "- CASH" (not "CASH AND CLEAR...").
$targetIndex = -1
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -eq "Amount Received mode`t`t- CASH`r") {
        $targetIndex = $i
        break
    }
}
Write-Host "Target paragraph index: $targetIndex"

if ($targetIndex -gt 0) {
    $anchor = $d.Paragraphs.Item($targetIndex).Range
    $anchor.Collapse(0)
    $anchor.InsertParagraphAfter()

    # Paragraph: blank separator line
    $p = $d.Paragraphs.Item($targetIndex + 1)
    $p.Range.InsertParagraphAfter()

    # Paragraph: date/time stamp line
    $p = $d.Paragraphs.Item($targetIndex + 2)
    $p.Range.Text = "TUE Feb 06"
    $p.Range.InsertAfter(" 11:11:18 PST 2018")
    $p.Range.InsertParagraphAfter()

    # Paragraph: Person Name line
    $p = $d.Paragraphs.Item($targetIndex + 3)
    $p.Range.Text = "Person Name"
    $p.Range.InsertAfter("`t")
    $p.Range.InsertAfter("`t")
    $p.Range.InsertAfter("`t")
    $p.Range.InsertAfter("`t")
    $p.Range.InsertAfter("- THA GOWRAMMA")
    $p.Range.InsertParagraphAfter()

    # Paragraph: dashed separator line
    $p = $d.Paragraphs.Item($targetIndex + 4)
    $p.Range.Text = "---------------------------------------------------------------"
    $p.Range.InsertParagraphAfter()

    # Paragraph: Item Name line
    $p = $d.Paragraphs.Item($targetIndex + 5)
    $p.Range.Text = "Item Name"
    $p.Range.InsertAfter("`t")
    $p.Range.InsertAfter("`t")
    $p.Range.InsertAfter("`t")
    $p.Range.InsertAfter("`t")
    $p.Range.InsertAfter("- CARROT EVE")
    $p.Range.InsertParagraphAfter()

    # Paragraph: Amount Received line (red text)
    $p = $d.Paragraphs.Item($targetIndex + 6)
    $p.Range.Font.Color = 255
    $p.Range.Text = "Amount Received"
    $p.Range.Font.Color = 255
    $p.Range.InsertAfter("`t")
    $p.Range.Font.Color = 255
    $p.Range.InsertAfter("`t")
    $p.Range.Font.Color = 255
    $p.Range.InsertAfter("`t")
    $p.Range.Font.Color = 255
    $p.Range.InsertAfter("- 556")
    $p.Range.InsertParagraphAfter()

    # Paragraph: Amount Received mode line
    $p = $d.Paragraphs.Item($targetIndex + 7)
    $p.Range.Text = "Amount Received mode"
    $p.Range.InsertAfter("`t")
    $p.Range.InsertAfter("`t")
    $p.Range.InsertAfter("- CASH AND CLEARD")
    $p.Range.InsertParagraphAfter()

    # Paragraph: trailing blank line
    $p = $d.Paragraphs.Item($targetIndex + 8)
    $p.Range.InsertParagraphAfter()

    Write-Host "New paragraph count: $($d.Paragraphs.Count)"
} else {
    Write-Host "ERROR: anchor paragraph not found"
}
